$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E6"   = 16.553
    "D7"   = -7.458
    "C9"   = -10.857
    "D12"  = -7.375999999999999
    "D14"  = -7.745000000000002
    "E15"  = 16.131
    "C18"  = -11.887
    "C20"  = -12.298
    "D26"  = -7.633
    "C27"  = -13.181
    "D27"  = -8.263
    "D29"  = -7.292
    "E33"  = 17.091
    "C35"  = -12.146
    "E35"  = 16.538
    "D37"  = -7.943000000000001
    "D38"  = -7.593999999999999
    "E38"  = 16.697
    "E43"  = 16.962
    "E44"  = 16.799
    "E47"  = 16.347
    "D51"  = -8.175000000000001
    "E51"  = 16.724
    "D52"  = -7.701000000000001
    "D55"  = -7.918000000000001
    "E57"  = 16.522
    "E63"  = 17.731
    "C69"  = -11.199
    "D69"  = -7.219999999999999
    "D70"  = -6.9
    "E70"  = 17.849
    "C76"  = -13.201
    "C78"  = -12.577
    "D81"  = -7.522
    "C82"  = -11.991
    "C83"  = -13.098
    "D83"  = -8.414
    "E88"  = 16.416
    "C93"  = -11.642
    "E99"  = 16.954
    "D102" = -7.662999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
